$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.546.84"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "1.568.84"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "212.59"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D8").Value = "46.02"
$ws.Range("E8").Value = "  +4.64%  "
$ws.Range("D9").Value = "24.07"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("D12").Value = "0.0887"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").Value = "1.792.66"
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("D14").Value = "1.593.06"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D16").Value = "28.526.92"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").Value = "230.25"
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("E21").Value = "  -2.63%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  -6.01%  "
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("E25").Value = "  +8.78%  "
$ws.Range("D26").Value = "151.29"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").Value = "6.42"
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("E29").Value = "  -3.38%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  +2.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.10"
$ws.Range("E32").Value = "  -3.07%  "
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.10"
$ws.Range("E34").Value = "  -1.47%  "
$ws.Range("D35").Value = "1.392.43"
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("E37").Value = "  -3.77%  "
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("E39").Value = "  +2.81%  "
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("E41").Value = "  -3.36%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").Value = "0.788"
$ws.Range("E45").Value = "  -4.03%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "0.969"
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("D48").Value = "62.81"
$ws.Range("E48").Value = "  -2.39%  "
$ws.Range("D49").Value = "1.705.53"
$ws.Range("E49").Value = "  -1.56%  "
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("E51").Value = "  -0.12%  "
